$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextCell 2 4 "26.772.98"
Set-TextCell 2 5 "  +0.12%  "
Set-TextCell 3 4 "1.538.39"
Set-TextCell 3 5 "  -1.77%  "
Set-TextCell 4 5 "  -0.02%  "
Set-TextCell 5 4 "205.64"
Set-TextCell 5 5 "  -0.33%  "
Set-TextCell 6 4 "0.485"
Set-TextCell 6 5 "  -0.72%  "
Set-TextCell 7 5 "  -0.02%  "
Set-TextCell 8 5 "  -0.37%  "
Set-TextCell 9 4 "21.24"
Set-TextCell 9 5 "  -2.73%  "
Set-TextCell 10 5 "  -0.54%  "
Set-TextCell 11 4 "0.0853"
Set-TextCell 11 5 "  -0.99%  "
Set-TextCell 12 4 "1.755.89"
Set-TextCell 12 5 "  -1.80%  "
Set-TextCell 13 4 "1.548.93"
Set-TextCell 13 5 "  -0.64%  "
Set-TextCell 14 5 "  -1.26%  "
Set-TextCell 15 5 "  -1.04%  "
Set-TextCell 16 4 "26.759.68"
Set-TextCell 16 5 "  -0.10%  "
Set-TextCell 17 4 "60.98"
Set-TextCell 17 5 "  -0.61%  "
Set-TextCell 18 4 "213.04"
Set-TextCell 18 5 "  -0.86%  "
Set-TextCell 19 4 "7.24"
Set-TextCell 19 5 "  -1.60%  "
Set-TextCell 20 4 "0.0$([char]0x2083)0682"
Set-TextCell 20 5 "  +0.85%  "
Set-TextCell 21 5 "  -0.01%  "
Set-TextCell 22 5 "  -2.00%  "
Set-TextCell 23 4 "9.15"
Set-TextCell 23 5 "  -1.55%  "
Set-TextCell 24 5 "  -3.22%  "
Set-TextCell 25 4 "151.54"
Set-TextCell 25 5 "  -0.63%  "
Set-TextCell 26 5 "  -2.08%  "
Set-TextCell 27 4 "14.77"
Set-TextCell 27 5 "  -0.97%  "
Set-TextCell 28 5 "  -0.03%  "
Set-TextCell 29 5 "  -0.91%  "
Set-TextCell 30 2 "PancakeSwap"
Set-TextCell 30 3 "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell 30 4 "1.10"
Set-TextCell 30 5 "  -0.90%  "
Set-TextCell 31 2 "Hedera"
Set-TextCell 31 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell 31 4 "0.0457"
Set-TextCell 31 5 "  -1.39%  "
Set-TextCell 32 5 "  +2.14%  "
Set-TextCell 33 4 "1.363.62"
Set-TextCell 33 5 "  -1.85%  "
Set-TextCell 34 5 "  +0.12%  "
Set-TextCell 35 5 "  -2.00%  "
Set-TextCell 36 5 "  +3.10%  "
Set-TextCell 37 5 "  -0.16%  "
Set-TextCell 38 5 "  +1.23%  "
Set-TextCell 39 5 "  -1.10%  "
Set-TextCell 40 4 "5.77"
Set-TextCell 40 5 "  +8.37%  "
Set-TextCell 41 5 "  -1.78%  "
Set-TextCell 42 4 "0.991"
Set-TextCell 42 5 "  +0.01%  "
Set-TextCell 43 4 "2.19"
Set-TextCell 43 5 "  +0.27%  "
Set-TextCell 44 4 "62.83"
Set-TextCell 44 5 "  -0.57%  "
Set-TextCell 45 4 "1.73"
Set-TextCell 45 5 "  -2.80%  "
Set-TextCell 46 4 "1.670.59"
Set-TextCell 46 5 "  -1.81%  "
Set-TextCell 47 4 "84.15"
Set-TextCell 47 5 "  -1.78%  "
Set-TextCell 48 4 "0.0509"
Set-TextCell 48 5 "  +3.58%  "
Set-TextCell 49 5 "  -1.40%  "
Set-TextCell 50 5 "  -0.59%  "
